$d = $word.ActiveDocument

# 1) Activation date update
$d.Content.Find.Execute("Ativação: 01/01/2018", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2022", 2)

# 2) Programa resumido (PT) - remove "Determinação de potência de bomba" item and renumber
$d.Content.Find.Execute("1)Determinação do número de Reynolds;2)Medidas de vazão em líquidos;3)Determinação de potência de bomba; 4)Esvaziamento de tanques cilíndricos;5)Determinação do fator de atrito em tubulações;6)Ensaios reológicos.", $true, $false, $false, $false, $false, $true, 1, $false, "1) Determinação do número de Reynolds; 2) Medidas de vazão em líquidos; 3) Esvaziamento de tanques cilíndricos; 4) Determinação do fator de atrito em tubulações; 5) Ensaios reológicos.", 2)

# 3) Programa resumido (EN, italic) - just add space after "1)"
$d.Content.Find.Execute("1)Determination of Reynolds number 2) Liquids flow measurements 3) Determination of the fluids velocity profile in pipe 4) Liquid flow from a cylindrical tank 5) Determination of friction factor in pipes 6) Rheological tests.", $true, $false, $false, $false, $false, $true, 1, $false, "1) Determination of Reynolds number 2) Liquids flow measurements 3) Determination of the fluids velocity profile in pipe 4) Liquid flow from a cylindrical tank 5) Determination of friction factor in pipes 6) Rheological tests.", 2)

# 4) Programa (PT) - remove "Determinação de potência de bomba centrífuga..." item and renumber
$d.Content.Find.Execute("1)Determinação do número de Reynolds, utilizando o aparato experimental;2)Medidas de Vazão em líquidos: determinação da vazão com utilização de placa de orifício, Venturi e rotâmetro;3)Determinação de potência de bomba centrífuga a partir da vazão volumétrica e das perdas de carga distribuída e localizada nas linhas de sucção e recalque de um sistema hidráulico;4)Esvaziamento de tanques cilíndricos: aplicação das equações de conservação de massa e energia na determinação do tempo de esvaziamento de tanques e comparação com dados experimentais;5)Determinação do fator de atrito em tubulações: avaliação do fator de atrito em função do número de Reynolds em tubulações;6)Ensaios reológicos: com a utilização de vários tipos de viscosímetros para determinar as viscosidades dinâmica e cinemática.", $true, $false, $false, $false, $false, $true, 1, $false, "1) Determinação do número de Reynolds, utilizando o aparato experimental; 2) Medidas de Vazão em líquidos: determinação da vazão com utilização de placa de orifício, Venturi e rotâmetro; 3) Esvaziamento de tanques cilíndricos: aplicação das equações de conservação de massa e energia na determinação do tempo de esvaziamento de tanques e comparação com dados experimentais; 4) Determinação do fator de atrito em tubulações: avaliação do fator de atrito em função do número de Reynolds em tubulações; 5) Ensaios reológicos: com a utilização de vários tipos de viscosímetros para determinar as viscosidades dinâmica e cinemática.", 2)

# 5) Programa (EN, italic) - remove "Determination of the fluids velocity profile in pipe..." item and renumber
$d.Content.Find.Execute("1) Determination of the Reynolds number using the experimental apparatus. 2) Liquids flow measurements: Determining the flow by using a rotameter, orifice plate and Venturi. 3) Determination of the fluids velocity profile in pipe: determining the velocity profile in pipe using a Pitot tube. 4) Liquid flow from a cylindrical tank: verification of a mathematical model based on the conservation equations of mass and energy to determine the time of emptying reservoir and compared the results with experimental data. 5) Determination of friction factor in pipes: evaluation of the friction factor as a function of Reynolds number in pipes. 6) Rheological tests: with the use of various types of viscometers to determine the dynamic and kinematic viscosities.", $true, $false, $false, $false, $false, $true, 1, $false, "1) Determination of the Reynolds number using the experimental apparatus. 2) Liquids flow measurements: Determining the flow by using a rotameter, orifice plate and Venturi. 3) Liquid flow from a cylindrical tank: verification of a mathematical model based on the conservation equations of mass and energy to determine the time of emptying reservoir and compared the results with experimental data. 4) Determination of friction factor in pipes: evaluation of the friction factor as a function of Reynolds number in pipes. 5) Rheological tests: with the use of various types of viscometers to determine the dynamic and kinematic viscosities.", 2)
